$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.819.09'
$ws.Range("E2").Value = '  -0.91%  '

# Row 3
$ws.Range("D3").Value = '2.351.05'
$ws.Range("E3").Value = '  -0.35%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.673'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.07'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.75%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.12%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.00'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.86%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.38%  '

# Row 15
$ws.Range("D15").Value = '2.700.73'
$ws.Range("E15").Value = '  -0.33%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.10%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.903'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.36%  '

# Row 18
$ws.Range("D18").Value = '2.350.51'
$ws.Range("E18").Value = '  -0.06%  '

# Row 19
$ws.Range("D19").Value = '43.753.07'
$ws.Range("E19").Value = '  -1.45%  '

# Row 20
$ws.Range("E20").Value = '  -1.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.96%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.64%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("E24").Value = '  +17.87%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.66%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.70%  '

# Row 28
$ws.Range("E28").Value = '  +5.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '

# Row 30
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '177.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.49%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '

# Row 32
$ws.Range("E32").Value = '  -1.46%  '

# Row 33
$ws.Range("E33").Value = '  +2.44%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.35%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.60%  '

# Row 37
$ws.Range("E37").Value = '  -2.18%  '

# Row 38
$ws.Range("E38").Value = '  -3.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.77%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0279'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.42%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '68.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +29.75%  '

# Row 42
$ws.Range("E42").Value = '  +11.68%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.36%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.01%  '

# Row 45
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.26%  '

# Row 46
$ws.Range("E46").Value = '  +3.34%  '

# Row 47
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.78%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.49'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.58%  '

# Row 49
$ws.Range("E49").Value = '  +0.05%  '

# Row 50
$ws.Range("E50").Value = '  -1.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.93%  '

Write-Host "Applied all changes"